$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C70").Value = 15731
$ws.Range("E70").Value = 24673307
$ws.Range("C79").Value = 116592
$ws.Range("E79").Value = 447357508
$ws.Range("C91").Value = 151119
$ws.Range("E91").Value = 482282585
$ws.Range("C92").Value = 409069
$ws.Range("D92").Value = 70906
$ws.Range("E92").Value = 1594718796
$ws.Range("C93").Value = 209545
$ws.Range("E93").Value = 1308331272
$ws.Range("C94").Value = 94175
$ws.Range("E94").Value = 916858080
$ws.Range("C95").Value = 50752
$ws.Range("E95").Value = 931855831
$ws.Range("C96").Value = 17256
$ws.Range("E96").Value = 790601490
$ws.Range("C97").Value = 2156
$ws.Range("E97").Value = 214088295
$ws.Range("C104").Value = 135232
$ws.Range("E104").Value = 272164522
$ws.Range("C114").Value = 3801
$ws.Range("E114").Value = 9113119
$ws.Range("C115").Value = 11693
$ws.Range("E115").Value = 32955791
$ws.Range("C165").Value = 83803
$ws.Range("D165").Value = 17113
$ws.Range("E165").Value = 354981515
$ws.Range("C167").Value = 12218
$ws.Range("E167").Value = 105743177
$ws.Range("C168").Value = 6205
$ws.Range("E168").Value = 100553127
